# "Atualizacao, removido Thread, adicionado scroll generico"
#
# Update the "Produtos" sheet: replace the retired "laptops" sample data
# row with a fresh one, and append a new validation row (14) used by the
# generic-scroll "no results" scenario. Also bump the stored username on
# the "Cadastro" sheet so the next automation run gets a unique login.

$wb = $excel.ActiveWorkbook

$produtos = $wb.Worksheets.Item("Produtos")

# Row 5 used to hold the old "HP ENVY x360" laptop fixture; swap it for
# the currently used "HP PAVILION 15T TOUCH LAPTOP" fixture (category,
# search term and expected result), matching the other LAPTOPS rows.
$produtos.Range("A5").Value = "LAPTOPS"
$produtos.Range("B5").Value = "HP PAVILION 15T TOUCH LAPTOP"
$produtos.Range("C5").Value = "HP PAVILION 15T TOUCH LAPTOP"

# Append a new row (14 / sheet row 15) for the generic "no results"
# scroll scenario.
$produtos.Cells.Item(15, 1).Value = "Geral"
$produtos.Cells.Item(15, 2).Value = "No results for"
$produtos.Cells.Item(15, 3).Value = "C(14,0) v(14,1)"

# Remember where the user was last working on this sheet, without
# stealing focus away from the "Cadastro" tab.
$produtos.Range("A12").Select()

$cadastro = $wb.Worksheets.Item("Cadastro")
$cadastro.Activate()

# Bump the stored username used by the cadastro (sign-up) test data.
$cadastro.Range("B2").Value = "Wilkerbn503"
